$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.903.07"
$ws.Range("E2").Value = "  -0.57%  "
$ws.Range("D3").Value = "2.301.24"
$ws.Range("E3").Value = "  -0.93%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "'299.76"
$ws.Range("E5").Value = "  -1.24%  "
$ws.Range("E6").Value = "  -2.34%  "
$ws.Range("E7").Value = "  +1.20%  "
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("D9").Value = "'0.507"
$ws.Range("E9").Value = "  -2.20%  "
$ws.Range("D10").Value = "'35.81"
$ws.Range("E10").Value = "  -1.36%  "
$ws.Range("E11").Value = "  -1.01%  "
$ws.Range("E12").Value = "  +0.71%  "
$ws.Range("D13").Value = "'17.66"
$ws.Range("E13").Value = "  -0.89%  "
$ws.Range("E14").Value = "  -2.31%  "
$ws.Range("D15").Value = "2.659.14"
$ws.Range("E15").Value = "  -0.91%  "
$ws.Range("D16").Value = "2.305.48"
$ws.Range("E16").Value = "  +0.19%  "
$ws.Range("D17").Value = "'0.778"
$ws.Range("E17").Value = "  -2.64%  "
$ws.Range("D18").Value = "42.874.88"
$ws.Range("E18").Value = "  -0.42%  "
$ws.Range("D19").Value = "'12.60"
$ws.Range("E19").Value = "  -4.10%  "
$ws.Range("D20").Value = "0.0₃0908"
$ws.Range("E20").Value = "  -0.58%  "
$ws.Range("D21").Value = "'6.10"
$ws.Range("E21").Value = "  -1.99%  "
$ws.Range("D22").Value = "'67.95"
$ws.Range("E22").Value = "  -0.40%  "
$ws.Range("D23").Value = "'241.94"
$ws.Range("E23").Value = "  +0.52%  "
$ws.Range("E24").Value = "  -1.51%  "
$ws.Range("E25").Value = "  +0.00%  "
$ws.Range("D26").Value = "'2.43"
$ws.Range("E26").Value = "  -1.28%  "
$ws.Range("E27").Value = "  -0.20%  "
$ws.Range("D28").Value = "'25.08"
$ws.Range("E28").Value = "  -1.71%  "
$ws.Range("D29").Value = "'165.96"
$ws.Range("E29").Value = "  -1.81%  "
$ws.Range("E30").Value = "  -0.91%  "
$ws.Range("E31").Value = "  -1.93%  "
$ws.Range("D32").Value = "'32.68"
$ws.Range("E32").Value = "  -4.75%  "
$ws.Range("E33").Value = "  +0.11%  "
$ws.Range("D34").Value = "'4.76"
$ws.Range("E34").Value = "  -3.89%  "
$ws.Range("E35").Value = "  -3.18%  "
$ws.Range("D36").Value = "'17.42"
$ws.Range("E36").Value = "  -2.57%  "
$ws.Range("D37").Value = "'2.39"
$ws.Range("E37").Value = "  -0.22%  "
$ws.Range("E38").Value = "  -1.91%  "
$ws.Range("E39").Value = "  -1.83%  "
$ws.Range("E40").Value = "  -3.72%  "
$ws.Range("E41").Value = "  -1.33%  "
$ws.Range("E42").Value = "  +0.25%  "
$ws.Range("D43").Value = "2.007.29"
$ws.Range("E43").Value = "  +0.57%  "
$ws.Range("E44").Value = "  -1.77%  "
$ws.Range("B45").Value = "FraxShare"
$ws.Range("C45").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D45").Value = "'10.18"
$ws.Range("E45").Value = "  -0.19%  "
$ws.Range("B46").Value = "ApeXProtocol"
$ws.Range("C46").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D46").Value = "'2.15"
$ws.Range("E46").Value = "  -4.30%  "
$ws.Range("D47").Value = "'17.15"
$ws.Range("E47").Value = "  -2.56%  "
$ws.Range("D48").Value = "'2.78"
$ws.Range("E48").Value = "  -2.98%  "
$ws.Range("D49").Value = "2.525.01"
$ws.Range("E49").Value = "  -0.93%  "
$ws.Range("D50").Value = "'53.34"
$ws.Range("E50").Value = "  -3.47%  "
$ws.Range("D51").Value = "'72.05"
$ws.Range("E51").Value = "  -5.59%  "
